# fix lỗi trong report cơ sở. Thêm cột ghi chú trong báo cáo về chi tiêu
#
# The workbook currently has 2 sheets: "Đơn sale chính", "Lương".
# After the edit it must have 3 sheets, in this order:
#   1) "Đơn sale chính"      (unchanged)
#   2) "Đơn phụ phẫu 1"      (new detail table -- takes over the old
#                             "Lương" sheet's position/identity)
#   3) "Lương"               (brand-new sheet, inserted after, holding the
#                             same payroll summary as before but with the
#                             SÓC TRĂNG figures now filled in)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: the existing "Lương" sheet becomes "Đơn phụ phẫu 1" and gets a
# brand-new detail table written into it (its old payroll content is
# relocated to a freshly-added "Lương" sheet in step 2).
# ---------------------------------------------------------------------
$donPhuPhau = $wb.Worksheets.Item(2)
$donPhuPhau.Cells.Clear()
$donPhuPhau.Name = "Đơn phụ phẫu 1"

$headers = @("Tiền tố", "Mã dịch vụ", "Ngày thực hiện", "Cơ sở", "Khách hàng", "Nguồn khách", "Tên dịch vụ", "Phụ phẫu 1", "Công phụ phẫu 1")
for ($col = 1; $col -le $headers.Length; $col++) {
    $donPhuPhau.Cells.Item(1, $col).Value = $headers[$col - 1]
}

# force the "Ngày thực hiện" (column C) data cells to be plain text so the
# dd-mm-yyyy strings below are kept literal instead of being auto-parsed
# into date serials (must be set on the exact target cells, not the whole
# column, for the format to actually take effect before the value write)
$donPhuPhau.Range("C2:C3").NumberFormat = "@"

$rows = @(
    @("HD-LUXURY", 626, "08-04-2024", "SÓC TRĂNG", "nguyễn thị mỹ trinh", "Cá nhân", "nhấn đồng tiền", "Trần Khánh Hiệp", $null),
    @("HD-LUXURY", 628, "08-04-2024", "SÓC TRĂNG", "nguyễn thị lệ trang", "Cá nhân", "Cắt mí", "Trần Khánh Hiệp", 50000)
)

$r = 2
foreach ($row in $rows) {
    for ($col = 1; $col -le $row.Length; $col++) {
        $val = $row[$col - 1]
        if ($null -ne $val) {
            $donPhuPhau.Cells.Item($r, $col).Value = $val
        }
    }
    $r++
}

# totals row
$donPhuPhau.Cells.Item(4, 1).Value = "Tổng"
$donPhuPhau.Cells.Item(4, 2).Value = 2
$donPhuPhau.Cells.Item(4, 9).Value = 50000

# ---------------------------------------------------------------------
# Step 2: insert a brand-new sheet right after "Đơn phụ phẫu 1", named
# "Lương", and repopulate it with the payroll summary (same as the
# original "Lương" sheet) with the SÓC TRĂNG section now computed.
# ---------------------------------------------------------------------
$luong = $wb.Worksheets.Add($null, $donPhuPhau)
$luong.Name = "Lương"

# match the sheetPr / pageMargins boilerplate carried by the other sheets
$luong.Outline.SummaryRow = 1
$luong.Outline.SummaryColumn = 1
$luong.PageSetup.LeftMargin = 54
$luong.PageSetup.RightMargin = 54
$luong.PageSetup.TopMargin = 72
$luong.PageSetup.BottomMargin = 72
$luong.PageSetup.HeaderMargin = 36
$luong.PageSetup.FooterMargin = 36

$luongRows = @(
    @(1, "Danh mục lương", 10),
    @(2, "Tổng công tại CẦN THƠ", 0),
    @(3, "Lương công tác tại CẦN THƠ", 0),
    @(4, "Lương cơ bản tại CẦN THƠ", $null),
    @(5, "Chiết khấu sale chính tại CẦN THƠ", 0),
    @(6, "Chiết khấu sale phụ tại CẦN THƠ", 0),
    @(7, "Đơn 1 bác sĩ tại CẦN THƠ", 0),
    @(8, "Đơn 2 bác sĩ tại CẦN THƠ", 0),
    @(9, "Công phụ phẫu 1 tại CẦN THƠ", 0),
    @(10, "Công phụ phẫu 2 tại CẦN THƠ", 0),
    @(11, "Ứng lương tại CẦN THƠ", 0),
    @(12, "Tổng công tại LONG XUYÊN", 0),
    @(13, "Lương công tác tại LONG XUYÊN", 0),
    @(14, "Lương cơ bản tại LONG XUYÊN", $null),
    @(15, "Chiết khấu sale chính tại LONG XUYÊN", 0),
    @(16, "Chiết khấu sale phụ tại LONG XUYÊN", 0),
    @(17, "Đơn 1 bác sĩ tại LONG XUYÊN", 0),
    @(18, "Đơn 2 bác sĩ tại LONG XUYÊN", 0),
    @(19, "Công phụ phẫu 1 tại LONG XUYÊN", 0),
    @(20, "Công phụ phẫu 2 tại LONG XUYÊN", 0),
    @(21, "Ứng lương tại LONG XUYÊN", 0),
    @(22, "Tổng công tại SÓC TRĂNG", 6),
    @(23, "Phụ cấp tại SÓC TRĂNG", 210000),
    @(24, "Lương cơ bản tại SÓC TRĂNG", 857142.8571428573),
    @(25, "Chiết khấu sale chính tại SÓC TRĂNG", 0),
    @(26, "Chiết khấu sale phụ tại SÓC TRĂNG", 0),
    @(27, "Đơn 1 bác sĩ tại SÓC TRĂNG", 0),
    @(28, "Đơn 2 bác sĩ tại SÓC TRĂNG", 0),
    @(29, "Công phụ phẫu 1 tại SÓC TRĂNG", 50000),
    @(30, "Công phụ phẫu 2 tại SÓC TRĂNG", 0),
    @(31, "Ứng lương tại SÓC TRĂNG", 0),
    @(32, "Tổng lương tại CẦN THƠ", 0),
    @(33, "Tổng lương tại LONG XUYÊN", 0),
    @(34, "Tổng lương tại SÓC TRĂNG", 1117142.857142857),
    @(35, "Tổng lương tại HỆ THỐNG", 1117142.857142857)
)

foreach ($entry in $luongRows) {
    $rowIdx = $entry[0]
    $label = $entry[1]
    $value = $entry[2]
    $luong.Cells.Item($rowIdx, 1).Value = $label
    if ($null -ne $value) {
        $luong.Cells.Item($rowIdx, 2).Value = $value
    }
}

# leave the first sheet selected, same as the source workbook
$wb.Worksheets.Item(1).Activate()
